$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row -> @{ D = new price text; E = new volume text }
# D values are prefixed with a leading apostrophe to force Excel to keep
# them as literal text (several look like plain numbers, e.g. "0.998") and
# the cell Style is reset to "Normal" right after so no stray
# quote-prefix / text-number-format style gets attached to the cell -
# matching the original inlineStr cells which carry no style index.
$updates = @(
    @{ Row = 2;  D = "78.364.36";  E = "  +2.60%  " },
    @{ Row = 3;  D = "3.168.14";   E = "  +6.60%  " },
    @{ Row = 4;  D = "0.998";      E = "  -0.21%  " },
    @{ Row = 5;  D = "204.19";     E = "  +2.32%  " },
    @{ Row = 6;  D = "633.29";     E = "  +0.88%  " },
    @{ Row = 7;  D = "0.998";      E = "  -0.13%  " },
    @{ Row = 8;  D = "0.228";      E = "  +14.17%  " },
    @{ Row = 9;  D = "0.575";      E = "  +5.10%  " },
    @{ Row = 10; D = "3.165.18";   E = "  +6.58%  " },
    @{ Row = 11; D = "0.573";      E = "  +33.40%  " },
    @{ Row = 12;                  E = "  +1.36%  " },
    @{ Row = 13; D = "5.45";       E = "  +9.85%  " },
    @{ Row = 14; D = "3.747.55";   E = "  +6.65%  " },
    @{ Row = 15; D = "0.0000227";  E = "  +21.41%  " },
    @{ Row = 16; D = "31.78";      E = "  +9.91%  " },
    @{ Row = 17; D = "78.142.77";  E = "  +2.41%  " },
    @{ Row = 18; D = "3.158.97";   E = "  +6.45%  " },
    @{ Row = 19; D = "14.32";      E = "  +7.08%  " },
    @{ Row = 20; D = "9.42";       E = "  +7.36%  " },
    @{ Row = 21; D = "430.66";     E = "  +16.09%  " },
    @{ Row = 22; D = "2.86";       E = "  +27.01%  " },
    @{ Row = 23; D = "4.87";       E = "  +13.67%  " },
    @{ Row = 24; D = "6.82";       E = "  +6.28%  " },
    @{ Row = 25; D = "3.327.30";   E = "  +6.57%  " },
    @{ Row = 26; D = "4.74";       E = "  +10.14%  " },
    @{ Row = 27; D = "76.41";      E = "  +5.10%  " },
    @{ Row = 28; D = "11.02";      E = "  +13.81%  " },
    @{ Row = 29;                  E = "  +0.18%  " },
    @{ Row = 30;                  E = "  +8.25%  " },
    @{ Row = 31; D = "0.998";      E = "  -0.42%  " },
    @{ Row = 32; D = "8.89";       E = "  +8.18%  " },
    @{ Row = 33; D = "1.49";       E = "  +7.78%  " },
    @{ Row = 34; D = "520.56";     E = "  +2.98%  " },
    @{ Row = 35; D = "1.98";       E = "  +2.95%  " },
    @{ Row = 36;                  E = "  +23.35%  " },
    @{ Row = 37; D = "22.52";      E = "  +11.12%  " },
    @{ Row = 38; D = "0.995";      E = "  -0.50%  " },
    @{ Row = 39; D = "0.397";      E = "  +4.78%  " },
    @{ Row = 40; D = "163.92";     E = "  +0.31%  " },
    @{ Row = 41; D = "196.66";     E = "  +6.40%  " },
    @{ Row = 42; D = "20.05";      E = "  +0.35%  " },
    @{ Row = 43;                  E = "  +4.15%  " },
    @{ Row = 45; D = "5.41";       E = "  +10.06%  " },
    @{ Row = 46; D = "0.800";      E = "  +14.29%  " },
    @{ Row = 47; D = "1.77";       E = "  +8.81%  " },
    @{ Row = 48; D = "1.30";       E = "  +6.08%  " },
    @{ Row = 49; D = "42.98";      E = "  +0.99%  " },
    @{ Row = 50; D = "2.57";       E = "  +11.60%  " },
    @{ Row = 51; D = "0.624";      E = "  +6.99%  " }
)

foreach ($u in $updates) {
    $row = $u.Row
    if ($u.ContainsKey("D")) {
        $cell = $ws.Range("D$row")
        $cell.Value = "'" + $u.D
        $cell.Style = "Normal"
    }
    if ($u.ContainsKey("E")) {
        $ws.Range("E$row").Value = $u.E
    }
}
